$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '65.913.55'
$r.Style = "Normal"
$r = $ws.Range('E2')
$r.NumberFormat = "@"
$r.Value = '  +6.54%  '
$r.Style = "Normal"
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.011.78'
$r.Style = "Normal"
$r = $ws.Range('E3')
$r.NumberFormat = "@"
$r.Value = '  +3.85%  '
$r.Style = "Normal"
$r = $ws.Range('E4')
$r.NumberFormat = "@"
$r.Value = '  +0.01%  '
$r.Style = "Normal"
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '583.64'
$r.Style = "Normal"
$r = $ws.Range('E5')
$r.NumberFormat = "@"
$r.Value = '  +2.86%  '
$r.Style = "Normal"
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '162.82'
$r.Style = "Normal"
$r = $ws.Range('E6')
$r.NumberFormat = "@"
$r.Value = '  +13.37%  '
$r.Style = "Normal"
$r = $ws.Range('E7')
$r.NumberFormat = "@"
$r.Value = '  -0.10%  '
$r.Style = "Normal"
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '3.007.41'
$r.Style = "Normal"
$r = $ws.Range('E8')
$r.NumberFormat = "@"
$r.Value = '  +3.80%  '
$r.Style = "Normal"
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.517'
$r.Style = "Normal"
$r = $ws.Range('E9')
$r.NumberFormat = "@"
$r.Value = '  +3.20%  '
$r.Style = "Normal"
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '6.90'
$r.Style = "Normal"
$r = $ws.Range('E10')
$r.NumberFormat = "@"
$r.Value = '  -0.38%  '
$r.Style = "Normal"
$r = $ws.Range('E11')
$r.NumberFormat = "@"
$r.Value = '  +7.95%  '
$r.Style = "Normal"
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.461'
$r.Style = "Normal"
$r = $ws.Range('E12')
$r.NumberFormat = "@"
$r.Value = '  +7.17%  '
$r.Style = "Normal"
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '0.0000252'
$r.Style = "Normal"
$r = $ws.Range('E13')
$r.NumberFormat = "@"
$r.Value = '  +9.22%  '
$r.Style = "Normal"
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '34.82'
$r.Style = "Normal"
$r = $ws.Range('E14')
$r.NumberFormat = "@"
$r.Value = '  +8.27%  '
$r.Style = "Normal"
$r = $ws.Range('E15')
$r.NumberFormat = "@"
$r.Value = '  -0.56%  '
$r.Style = "Normal"
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '65.922.90'
$r.Style = "Normal"
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '6.98'
$r.Style = "Normal"
$r = $ws.Range('E18')
$r.NumberFormat = "@"
$r.Value = '  +7.10%  '
$r.Style = "Normal"
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '3.009.33'
$r.Style = "Normal"
$r = $ws.Range('E19')
$r.NumberFormat = "@"
$r.Value = '  +4.53%  '
$r.Style = "Normal"
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '457.53'
$r.Style = "Normal"
$r = $ws.Range('E20')
$r.NumberFormat = "@"
$r.Value = '  +6.21%  '
$r.Style = "Normal"
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '13.95'
$r.Style = "Normal"
$r = $ws.Range('E21')
$r.NumberFormat = "@"
$r.Value = '  +8.10%  '
$r.Style = "Normal"
$r = $ws.Range('E22')
$r.NumberFormat = "@"
$r.Value = '  +5.52%  '
$r.Style = "Normal"
$r = $ws.Range('E23')
$r.NumberFormat = "@"
$r.Value = '  +7.56%  '
$r.Style = "Normal"
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '82.55'
$r.Style = "Normal"
$r = $ws.Range('E24')
$r.NumberFormat = "@"
$r.Value = '  +4.62%  '
$r.Style = "Normal"
$r = $ws.Range('E25')
$r.NumberFormat = "@"
$r.Value = '  +15.35%  '
$r.Style = "Normal"
$r = $ws.Range('E26')
$r.NumberFormat = "@"
$r.Value = '  +3.29%  '
$r.Style = "Normal"
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '10.60'
$r.Style = "Normal"
$r = $ws.Range('E27')
$r.NumberFormat = "@"
$r.Value = '  +5.55%  '
$r.Style = "Normal"
$r = $ws.Range('E28')
$r.NumberFormat = "@"
$r.Value = '  -0.10%  '
$r.Style = "Normal"
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '8.13'
$r.Style = "Normal"
$r = $ws.Range('E29')
$r.NumberFormat = "@"
$r.Value = '  +16.11%  '
$r.Style = "Normal"
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '2.33'
$r.Style = "Normal"
$r = $ws.Range('E30')
$r.NumberFormat = "@"
$r.Value = '  +15.67%  '
$r.Style = "Normal"
$r = $ws.Range('E31')
$r.NumberFormat = "@"
$r.Value = '  +4.22%  '
$r.Style = "Normal"
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '0.0000102'
$r.Style = "Normal"
$r = $ws.Range('E32')
$r.NumberFormat = "@"
$r.Value = '  -7.56%  '
$r.Style = "Normal"
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '27.08'
$r.Style = "Normal"
$r = $ws.Range('E33')
$r.NumberFormat = "@"
$r.Value = '  +5.76%  '
$r.Style = "Normal"
$r = $ws.Range('E34')
$r.NumberFormat = "@"
$r.Value = '  +3.62%  '
$r.Style = "Normal"
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$r = $ws.Range('E35')
$r.NumberFormat = "@"
$r.Value = '  -0.05%  '
$r.Style = "Normal"
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '0.994'
$r.Style = "Normal"
$r = $ws.Range('E36')
$r.NumberFormat = "@"
$r.Value = '  +4.16%  '
$r.Style = "Normal"
$r = $ws.Range('E37')
$r.NumberFormat = "@"
$r.Value = '  +7.62%  '
$r.Style = "Normal"
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '2.15'
$r.Style = "Normal"
$r = $ws.Range('E38')
$r.NumberFormat = "@"
$r.Value = '  +12.61%  '
$r.Style = "Normal"
$r = $ws.Range('E39')
$r.NumberFormat = "@"
$r.Value = '  +6.03%  '
$r.Style = "Normal"
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '49.83'
$r.Style = "Normal"
$r = $ws.Range('E40')
$r.NumberFormat = "@"
$r.Value = '  +2.07%  '
$r.Style = "Normal"
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '0.310'
$r.Style = "Normal"
$r = $ws.Range('E41')
$r.NumberFormat = "@"
$r.Value = '  +15.87%  '
$r.Style = "Normal"
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '0.121'
$r.Style = "Normal"
$r = $ws.Range('E42')
$r.NumberFormat = "@"
$r.Value = '  +5.95%  '
$r.Style = "Normal"
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '43.64'
$r.Style = "Normal"
$r = $ws.Range('E43')
$r.NumberFormat = "@"
$r.Value = '  +8.23%  '
$r.Style = "Normal"
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '8.49'
$r.Style = "Normal"
$r = $ws.Range('E44')
$r.NumberFormat = "@"
$r.Value = '  +4.25%  '
$r.Style = "Normal"
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '387.71'
$r.Style = "Normal"
$r = $ws.Range('E45')
$r.NumberFormat = "@"
$r.Value = '  +12.32%  '
$r.Style = "Normal"
$r = $ws.Range('E46')
$r.NumberFormat = "@"
$r.Value = '  +6.26%  '
$r.Style = "Normal"
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '2.794.90'
$r.Style = "Normal"
$r = $ws.Range('E47')
$r.NumberFormat = "@"
$r.Value = '  +3.43%  '
$r.Style = "Normal"
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '135.18'
$r.Style = "Normal"
$r = $ws.Range('E48')
$r.NumberFormat = "@"
$r.Value = '  +2.69%  '
$r.Style = "Normal"
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '24.01'
$r.Style = "Normal"
$r = $ws.Range('E50')
$r.NumberFormat = "@"
$r.Value = '  +11.40%  '
$r.Style = "Normal"
$r = $ws.Range('E51')
$r.NumberFormat = "@"
$r.Value = '  +3.96%  '
$r.Style = "Normal"
